$d = $word.ActiveDocument

# Append three new paragraphs after the existing "SerialNumber" paragraph:
#   - "ConfName"
#   - "ConfDate"
#   - a trailing empty paragraph
# All inherit the fr-FR language formatting already used in the document.
# Using Find/Replace with "^p" paragraph-mark codes (instead of separate
# InsertParagraphAfter calls) so the final empty paragraph is created
# without a leftover empty run.
$d.Content.Find.Execute("SerialNumber", $true, $false, $false, $false, $false,
                         $true, 1, $false, "SerialNumber^pConfName^pConfDate^p", 2)
